# WhitePack_Full_Template_v7.2.xlsx — "Add files via upload" re-upload edit
#
# Content changes applied (per the OOXML diff):
#  1. Products sheet, cell J2: the shared "image URL" text/hyperlink display
#     text is swapped for a new image link (the underlying hyperlink target
#     itself is untouched, only the visible cell text / shared string changes).
#  2. Products sheet: the saved cursor/selection moves from J6 to C2.
#  3. Products sheet: column J (10) is widened, to comfortably show the new,
#     longer URL text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# 1) Update the displayed text of the hyperlinked image-URL cell.
$ws.Range("J2").Value = "https://sanitalb.com/public/uploads/images/69428894249722360.jpg"

# 3) Widen column J so the longer URL text fits better.
$ws.Columns.Item(10).ColumnWidth = 74.25

# 2) Move/save the active selection to C2 (this also matches where the sheet
#    is left selected when the file is re-saved).
$ws.Range("C2").Select()
